$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "dfg"
$ws.Range("D2").Value = "dfg"

# Update row 3
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "3"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "updated"
$ws.Range("D3").Value = "dfg"

# Update row 4
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "4"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "updated"
$ws.Range("D4").Value = "dfg"

# Add new row 5
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 4
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "5"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "dfg"
$ws.Range("D5").Value = "dfg"
